$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 'Última actualización: 07:50:27'
$ws.Range("A3").Value = 'Total filas: 94'
$data1 = @(
    @('04:03:00','04:03','81_EL PELIGRO',0,'LP1912'),
    @('04:37:19','04:46','215A_EL PATO',9,'LP1912'),
    @('04:52:25','04:53','11_ETCHEVERRY',1,'LP1912'),
    @('04:52:25','05:16','17_ROMERO',24,'LP1912'),
    @('05:20:30','05:20','17_ROMERO',0,'LP1912'),
    @('04:52:25','05:22','23_HERNANDEZ',30,'LP1912'),
    @('05:20:30','05:26','23_HERNANDEZ',6,'LP1912'),
    @('05:20:30','05:34','215B_EL PATO',14,'LP1912'),
    @('04:03:00','05:35','215B_EL PATO',92,'LP1912'),
    @('04:03:00','05:41','14_ABASTO',98,'LP1912'),
    @('05:20:30','05:46','15_ABASTO',26,'LP1912'),
    @('05:54:55','05:54','10_OLMOS',0,'LP1912'),
    @('05:54:55','05:55','81_EL PELIGRO',1,'LP1912'),
    @('05:20:30','06:04','16_SANTA ANA',44,'LP1912'),
    @('05:54:55','06:11','215A_EL PATO',17,'LP1912'),
    @('05:54:55','06:13','225_HARAS DEL SUR',19,'LP1912'),
    @('05:20:30','06:14','225_HARAS DEL SUR',54,'LP1912'),
    @('05:54:55','06:20','26_HERNANDEZ',26,'LP1912'),
    @('05:20:30','06:21','26_HERNANDEZ',61,'LP1912'),
    @('05:54:55','06:26','23_HERNANDEZ',32,'LP1912'),
    @('06:24:16','06:27','23_HERNANDEZ',3,'LP1912'),
    @('05:54:55','06:29','86_EST CHICA-ESC AGRARIA',35,'LP1912'),
    @('06:24:16','06:30','86_EST CHICA-ESC AGRARIA',6,'LP1912'),
    @('06:24:16','06:31','16_SANTA ANA',7,'LP1912'),
    @('05:54:55','06:43','225_C ROCA-H SUR',49,'LP1912'),
    @('06:24:16','06:44','225_C ROCA-H SUR',20,'LP1912'),
    @('06:24:16','06:46','215C_EL PATO',22,'LP1912'),
    @('06:53:31','06:56','14_ABASTO',3,'LP1912'),
    @('05:54:55','06:59','14_ABASTO',65,'LP1912'),
    @('06:24:16','07:00','14_ABASTO',36,'LP1912'),
    @('06:53:31','07:01','16_SANTA ANA',8,'LP1912'),
    @('06:53:31','07:04','23_HERNANDEZ',11,'LP1912'),
    @('06:53:31','07:05','15_ABASTO',12,'LP1912'),
    @('05:54:55','07:06','225_GOMEZ',72,'LP1912'),
    @('06:53:31','07:07','225_GOMEZ',14,'LP1912'),
    @('06:53:31','07:11','215A_EL PATO',18,'LP1912'),
    @('06:53:31','07:15','11_ETCHEVERRY',22,'LP1912'),
    @('06:53:31','07:16','16_SANTA ANA',23,'LP1912'),
    @('06:24:16','07:16','11_ETCHEVERRY',52,'LP1912'),
    @('07:18:07','07:18','16_SANTA ANA',0,'LP1912'),
    @('05:54:55','07:20','26_HERNANDEZ',86,'LP1912'),
    @('07:18:07','07:21','26_HERNANDEZ',3,'LP1912'),
    @('07:18:07','07:21','10_OLMOS',3,'LP1912'),
    @('07:18:07','07:21','225_GOMEZ',3,'LP1912'),
    @('05:54:55','07:22','10_OLMOS',88,'LP1912'),
    @('06:53:31','07:23','10_OLMOS',30,'LP1912'),
    @('05:54:55','07:31','16_SANTA ANA',97,'LP1912'),
    @('06:53:31','07:31','11_ETCHEVERRY',38,'LP1912'),
    @('07:18:07','07:32','16_SANTA ANA',14,'LP1912'),
    @('07:18:07','07:32','11_ETCHEVERRY',14,'LP1912'),
    @('07:18:07','07:32','84_COLONIA URQUIZA-ESC 49',14,'LP1912'),
    @('06:53:31','07:36','27_EL RETIRO',43,'LP1912'),
    @('07:18:07','07:36','23_HERNANDEZ',18,'LP1912'),
    @('07:18:07','07:37','27_EL RETIRO',19,'LP1912'),
    @('06:24:16','07:39','10_OLMOS',75,'LP1912'),
    @('06:53:31','07:47','14_ABASTO',54,'LP1912'),
    @('07:18:07','07:48','14_ABASTO',30,'LP1912'),
    @('07:50:27','07:50','16_SANTA ANA',0,'LP1912'),
    @('07:50:27','07:51','215D_EL PATO',1,'LP1912'),
    @('07:18:07','07:52','215D_EL PATO',34,'LP1912'),
    @('07:50:27','07:52','10_OLMOS',2,'LP1912'),
    @('07:18:07','07:55','10_OLMOS',37,'LP1912'),
    @('07:50:27','07:59','23_HERNANDEZ',9,'LP1912'),
    @('07:18:07','08:00','23_HERNANDEZ',42,'LP1912'),
    @('07:50:27','08:03','11_ETCHEVERRY',13,'LP1912'),
    @('07:18:07','08:04','11_ETCHEVERRY',46,'LP1912'),
    @('06:53:31','08:05','23_HERNANDEZ',72,'LP1912'),
    @('07:50:27','08:10','16_SANTA ANA',20,'LP1912'),
    @('07:50:27','08:12','15_ABASTO',22,'LP1912'),
    @('07:50:27','08:13','10_OLMOS',23,'LP1912'),
    @('07:50:27','08:21','26_HERNANDEZ',31,'LP1912'),
    @('07:50:27','08:22','16_P MOR-SANTA ANA',32,'LP1912'),
    @('07:18:07','08:23','16_P MOR-SANTA ANA',65,'LP1912'),
    @('07:50:27','08:23','215B_EL PATO',33,'LP1912'),
    @('07:50:27','08:27','84_COLONIA URQUIZA-ESC 49',37,'LP1912'),
    @('07:50:27','08:35','23_HERNANDEZ',45,'LP1912'),
    @('07:18:07','08:42','81_EL PELIGRO',84,'LP1912'),
    @('07:50:27','08:43','14_ABASTO',53,'LP1912'),
    @('07:18:07','08:44','14_ABASTO',86,'LP1912'),
    @('07:50:27','08:54','17_ROMERO',64,'LP1912'),
    @('07:50:27','09:01','215A_EL PATO',71,'LP1912'),
    @('07:18:07','09:02','215A_EL PATO',104,'LP1912'),
    @('07:50:27','09:10','16_P MOR-SANTA ANA',80,'LP1912'),
    @('07:18:07','09:11','16_P MOR-SANTA ANA',113,'LP1912'),
    @('07:50:27','09:11','81_EL PELIGRO',81,'LP1912'),
    @('07:50:27','09:16','27_EL RETIRO',86,'LP1912'),
    @('07:18:07','09:17','27_EL RETIRO',119,'LP1912'),
    @('07:50:27','09:21','26_HERNANDEZ',91,'LP1912'),
    @('07:50:27','09:22','17_ROMERO',92,'LP1912'),
    @('07:50:27','09:23','11_ETCHEVERRY',93,'LP1912'),
    @('07:50:27','09:32','15_ABASTO',102,'LP1912'),
    @('07:50:27','09:33','10_OLMOS',103,'LP1912'),
    @('07:50:27','09:42','215C_EL PATO',112,'LP1912'),
    @('07:50:27','09:43','14_ABASTO',113,'LP1912')
)
for ($i = 0; $i -lt $data1.Length; $i++) {
    $row = $data1[$i]
    $r = 6 + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# ---- Sheet 2: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 'Última actualización: 07:50:27'
$ws.Range("A3").Value = 'Total filas: 12'
$data2 = @(
    @('04:37:19','04:46','215A_EL PATO',9,'LP1912'),
    @('05:20:30','05:34','215B_EL PATO',14,'LP1912'),
    @('04:03:00','05:35','215B_EL PATO',92,'LP1912'),
    @('05:54:55','06:11','215A_EL PATO',17,'LP1912'),
    @('06:24:16','06:46','215C_EL PATO',22,'LP1912'),
    @('06:53:31','07:11','215A_EL PATO',18,'LP1912'),
    @('07:50:27','07:51','215D_EL PATO',1,'LP1912'),
    @('07:18:07','07:52','215D_EL PATO',34,'LP1912'),
    @('07:50:27','08:23','215B_EL PATO',33,'LP1912'),
    @('07:50:27','09:01','215A_EL PATO',71,'LP1912'),
    @('07:18:07','09:02','215A_EL PATO',104,'LP1912'),
    @('07:50:27','09:42','215C_EL PATO',112,'LP1912')
)
for ($i = 0; $i -lt $data2.Length; $i++) {
    $row = $data2[$i]
    $r = 6 + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# ---- Sheet 3: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 'Última actualización: 07:50:27'
$ws.Range("A3").Value = 'Total filas: 14'
$data3 = @(
    @('05:20:30','05:44','215A_LA PLATA',24,'L6173'),
    @('05:54:55','06:09','215A_LA PLATA',15,'L6173'),
    @('05:20:30','06:10','215A_LA PLATA',50,'L6173'),
    @('05:54:55','06:32','215C_LA PLATA',38,'L6203'),
    @('06:24:16','06:33','215C_LA PLATA',9,'L6203'),
    @('05:54:55','06:59','215B_LP-P MOR-1 Y 57',65,'L6173'),
    @('06:53:31','07:00','215B_LP-P MOR-1 Y 57',7,'L6173'),
    @('05:54:55','07:34','215A_LA PLATA',100,'L6173'),
    @('07:18:07','07:35','215A_LA PLATA',17,'L6173'),
    @('06:53:31','08:07','215C_LA PLATA',74,'L6203'),
    @('07:50:27','08:11','215C_LA PLATA',21,'L6203'),
    @('07:18:07','08:14','215C_LA PLATA',56,'L6203'),
    @('07:50:27','08:35','215A_LA PLATA',45,'L6173'),
    @('07:50:27','09:09','215D_LA PLATA',79,'L6203')
)
for ($i = 0; $i -lt $data3.Length; $i++) {
    $row = $data3[$i]
    $r = 6 + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

